# Clava Parser Restructuring — convert SwitchStmt, CaseStmt, DefaultStmt (and
# SwitchCase) in the "stmt" sheet, plus SizeOfPackExpr in the "expr" sheet,
# from pending to completed ("o"), and rename the "expr" sheet to "-expr-"
# to match the other "-name-" style tabs.

$wb = $excel.ActiveWorkbook

# 1. Rename the "expr" sheet to "-expr-". Excel automatically rewrites any
#    formulas that reference the sheet by name (e.g. the Summary totals).
$exprSheet = $wb.Worksheets.Item("expr")
$exprSheet.Name = "-expr-"

# 2. Mark the newly-converted nodes as completed ("o") on the "-expr-" sheet.
$exprSheet.Range("B59").Value = "o"

# 3. Mark the newly-converted nodes as completed ("o") on the "stmt" sheet:
#    row 4 = CaseStmt, row 12 = DefaultStmt, row 24 = SwitchCase,
#    row 25 = SwitchStmt.
$stmtSheet = $wb.Worksheets.Item("stmt")
$stmtSheet.Range("B4").Value = "o"
$stmtSheet.Range("B12").Value = "o"
$stmtSheet.Range("B24").Value = "o"
$stmtSheet.Range("B25").Value = "o"

# 4. Update the selection left behind on the "-expr-" sheet (no longer the
#    active tab) and make "stmt" the active sheet/selection, matching where
#    the author was working when they saved.
$exprSheet.Activate()
$exprSheet.Range("B44").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$stmtSheet.Activate()
$stmtSheet.Range("B26").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
